# "Add files via upload" — the tracked OOXML delta for this commit is a
# deletion of the second picture ("圖片 4" / id=5, rId3) that was added to
# slide 5 ("Exploratory Data Analysis (EDA)"). Everything else in the
# canonical diff (ma14:wrappingTextBoxFlag attribute-order churn, a stray
# endParaRPr dirty="0") is PowerPoint's own re-serialization noise from the
# save that produced this commit, not a deliberate edit, so we only need to
# remove that shape here.

$p = $ppt.ActivePresentation

# Slide 5 = "Exploratory Data Analysis (EDA)" (rels: rId6 -> slides/slide5.xml,
# which is the 5th entry in p:sldIdLst).
$s = $p.Slides.Item(5)

# Walk the shapes and remove the picture named "圖片 4" (creationId
# {AB698198-6B17-A182-54F6-672E24699659}, blip r:embed="rId3" ->
# ../media/image7.png) rather than hard-coding an index, so the script is
# robust to minor shape-ordering differences.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "圖片 4") {
        $shape.Delete()
    }
}
